# Zhihong_Deng_Handover.xlsx edit:
#  - Column I ("Assigned to which Continuing Member?") is repurposed to record
#    the person who is *passing* the inheritance rather than the one
#    inheriting it, and every data row is updated from "Muhammad Arslan" to
#    "Zhihong Deng" (Mingbo's inheritance handover, per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for column I (row 3)
$ws.Range("I3").Value = "Name of the person passing the inheritance."

# Data rows 4-12: replace the assigned person's name
$ws.Range("I4").Value = "Zhihong Deng"
$ws.Range("I5").Value = "Zhihong Deng"
$ws.Range("I6").Value = "Zhihong Deng"
$ws.Range("I7").Value = "Zhihong Deng"
$ws.Range("I8").Value = "Zhihong Deng"
$ws.Range("I9").Value = "Zhihong Deng"
$ws.Range("I10").Value = "Zhihong Deng"
$ws.Range("I11").Value = "Zhihong Deng"
$ws.Range("I12").Value = "Zhihong Deng"

# Update on-screen selection/scroll position to match the author's final
# cursor location when they saved the workbook.
$ws.Range("I11").Select() | Out-Null
